# Update with author-specific affiliations
# Column C ("Faculty") is replaced with "University of North Carolina"
# for every data row, and the active selection moves to F11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newAffiliation = "University of North Carolina"

for ($row = 1; $row -le 22; $row++) {
    $ws.Cells.Item($row, 3).Value = $newAffiliation
}

$ws.Range("F11").Select()
